$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: replace the old "placeholder" shared-formula row with real figures ---
$ws.Range("A10").Value = 1260.6300000000001
$ws.Range("B10").Value = 835.31
$ws.Range("C10").Formula = "=A10+B10"
$ws.Range("E10").Value = 1007.96
$ws.Range("F10").Value = 835.31
$ws.Range("G10").Formula = "=E10+F10"

# --- Row 11: new formulas replacing the old shared placeholder ---
$ws.Range("A11").Formula = "=((1248.94+(2579.69/5))/5)*8"
$ws.Range("B11").Formula = "=B6/5*8"
$ws.Range("C11").Formula = "=A11+B11"
$ws.Range("E11").Value = 2708.94
$ws.Range("F11").Value = 1871.1
$ws.Range("G11").Formula = "=E11+F11"

# --- Row 12: totals row - B12 was empty, now sums the Mitarbeiter column ---
$ws.Range("B12").Formula = "=SUM(B3:B11)"

# --- Row 14: new "Gesamt-Budget + 10%" style delta cell next to the Projektzeit label ---
$ws.Range("D14").Formula = "=A11-E11"
$ws.Range("D14").NumberFormat = '#,##0.00\ "' + [char]0x20AC + '"'

# --- New row 30: quick scratch calculation added at the bottom of the sheet ---
$ws.Range("E30").Formula = "=22476.38-22571.22"

# --- Leave the cursor where the author left it when they saved the file ---
[void]$ws.Range("E31").Select()
